$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column A with value 1 for rows 2-4, 6-9, 11-14 (skipping 5 and 10)
$rows = @(2,3,4,6,7,8,9,11,12,13,14)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = 1
}

# Formula in D8 counting occurrences of A4's value within A2:A14
$ws.Range("D8").Formula = "=COUNTIF(A2:A14,A4)"

# Selection as captured in the saved view state
$ws.Range("D9").Select()
